$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "Image"
$ws.Range("I1").Value = "DOB"

$ws.Range("C2").Value = "ash"
$ws.Range("D2").Value = "rai"
$ws.Range("E2").Value = "bac"
$ws.Range("F2").Value = "PhD"
$ws.Range("H2").Value = "iii"
$ws.Range("I2").Value = "28-04-2028"

$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Select()
